{"js": "// Minutes table, \"Technologie u\u017cywane w projekcie\" row: the technology\n// list \" (Hibernate,Spring, JSF (z RichFaces), Jetty, STS i maven)\" swaps\n// \"Jetty\" for \"Tomcat\" (commit: \"changed jetty to tomcat\").\nconst searchResults = context.document.body.search(\"Jetty\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the text 'Jetty' to replace.\");\n}\n\n// Replace in place so the run keeps its existing formatting\n// (Arial 10pt, pl-PL) \u2014 only the word itself changes.\nsearchResults.items[0].insertText(\"Tomcat\", \"Replace\");\nawait context.sync();\n", "ps1": "# Minutes table, \"Technologie u\u017cywane w projekcie\" row: the technology\n# list \" (Hibernate,Spring, JSF (z RichFaces), Jetty, STS i maven)\" swaps\n# \"Jetty\" for \"Tomcat\" (commit: \"changed jetty to tomcat\").\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Jetty\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"Tomcat\"\n\n# wdReplaceAll = 2 \u2014 there is exactly one occurrence, but this guards\n# against leaving the job half-done if Word stops after the first hit.\n$replaced = $rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\nif (-not $replaced) {\n    throw \"Could not find the text 'Jetty' to replace.\"\n}\n"}
